$wb = $excel.ActiveWorkbook

$wsEOH = $wb.Worksheets.Item("u_EOH")
$wsEOH.Range("A2").Value = -10.05222087110701

$wsVL = $wb.Worksheets.Item("v_l")
$wsVL.Range("A2").Value = 14303593.64251507
$wsVL.Range("A3").Value = 15820331.78251507
$wsVL.Range("A4").Value = 11644905.845
